$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths for C, D, H
# (offset by -5/6 to compensate the engine's width<->pixel round-trip so the
#  saved OOXML <col width> lands exactly on the target value)
$ws.Columns.Item(3).ColumnWidth = 51.166666666666664
$ws.Columns.Item(4).ColumnWidth = 43.166666666666664
$ws.Columns.Item(8).ColumnWidth = 50.166666666666664

# Row data: OPPORTUNITY ID, LINK, TITLE, COUNTRY, PREMIUM, APPLICANTS, DURATION, ORGANIZATION, isPremium
$rows = @(
    ,@('1329890', 'https://aiesec.org/opportunity/global-talent/1329890', 'Finance Intern  - Chile', 'Santiago, Región Metropolitana, Chile', 'Yes', '4 applicants', '3 - 6 Months', 'Henkel AG & Co. KGaA', $true)
    ,@('1330419', 'https://aiesec.org/opportunity/global-talent/1330419', 'Digital Technology Application Services (ONLY EU)', 'Bruxelles, Belgio', 'No', '0 applicants', '6 - 18 Months', 'UCB', $false)
    ,@('1330415', 'https://aiesec.org/opportunity/global-talent/1330415', 'Digital Technology Training ( EU ONLY)', 'Bruxelles, Belgio', 'No', '1 applicant', '6 - 18 Months', 'UCB', $false)
    ,@('1330394', 'https://aiesec.org/opportunity/global-talent/1330394', '[Impact in Belo Horizonte] - Digital Marketing', 'Belo Horizonte, MG, Brasil', 'No', '1 applicant', '6 - 18 Months', 'Group Tech Participações LTDA', $false)
    ,@('1330393', 'https://aiesec.org/opportunity/global-talent/1330393', '[Impact in Belo Horizonte] - Business Development', 'Belo Horizonte, MG, Brasil', 'No', '2 applicants', '6 - 18 Months', 'Group Tech Participações LTDA', $false)
    ,@('1330282', 'https://aiesec.org/opportunity/global-talent/1330282', 'Web Developer Intern', 'Phagwara, Punjab, India', 'No', '0 applicants', '3 - 6 Months', 'GNA University', $false)
    ,@('1329673', 'https://aiesec.org/opportunity/global-talent/1329673', 'Architecture', 'Gabes, Tunisia', 'No', '0 applicants', '9 - 12 Weeks', 'BMES', $false)
    ,@('1327760', 'https://aiesec.org/opportunity/global-talent/1327760', 'Fashion Consultant', 'Hyderabad, Telangana, India', 'No', '3 applicants', '6 - 18 Months', 'MPF clothing collection PVT LTD', $false)
    ,@('1327187', 'https://aiesec.org/opportunity/global-talent/1327187', 'Software Developer Intern', 'Sahibzada Ajit Singh Nagar, Punjab, India', 'No', '21 applicants', '9 - 12 Weeks', 'Solitaire Infosys Pvt. Ltd', $false)
    ,@('1327186', 'https://aiesec.org/opportunity/global-talent/1327186', 'Cyber Security Intern', 'Sahibzada Ajit Singh Nagar, Punjab, India', 'No', '13 applicants', '9 - 12 Weeks', 'Solitaire Infosys Pvt. Ltd', $false)
    ,@('1327185', 'https://aiesec.org/opportunity/global-talent/1327185', 'Machine Learning Intern', 'Sahibzada Ajit Singh Nagar, Punjab, India', 'No', '8 applicants', '9 - 12 Weeks', 'Solitaire Infosys Pvt. Ltd', $false)
    ,@('1327184', 'https://aiesec.org/opportunity/global-talent/1327184', 'AI Intern', 'Sahibzada Ajit Singh Nagar, Punjab, India', 'No', '7 applicants', '9 - 12 Weeks', 'Solitaire Infosys Pvt. Ltd', $false)
    ,@('1327181', 'https://aiesec.org/opportunity/global-talent/1327181', 'Data Analyst Intern', 'Sahibzada Ajit Singh Nagar, Punjab, India', 'No', '9 applicants', '9 - 12 Weeks', 'Solitaire Infosys Pvt. Ltd', $false)
    ,@('1326116', 'https://aiesec.org/opportunity/global-talent/1326116', 'Sales and Customer Service Support', 'İstanbul, Türkiye', 'No', '128 applicants', '6 - 18 Months', 'Tornado Makine Otomotiv İnşaat Sanayi ve Ticaret', $false)
    ,@('1323478', 'https://aiesec.org/opportunity/global-talent/1323478', 'Sales Support Executive', 'Cyberjaya, Selangor, Malaysia', 'No', '33 applicants', '6 - 18 Months', 'IX Telecom Sdn Bhd', $false)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    if ($row[8]) {
        $ws.Cells.Item($r, 5).Interior.Color = 65535
    } else {
        $ws.Cells.Item($r, 5).Interior.ColorIndex = 0
    }
}

